$wb = $excel.ActiveWorkbook

# The "Poland" sheet was selected (Ctrl+A / Select All) and copied to create
# the new "UK" sheet, which was placed immediately after it (as the last tab).
$poland = $wb.Worksheets.Item("Poland")
$poland.Copy($null, $poland)

$uk = $wb.Worksheets.Item("Poland (2)")
$uk.Name = "UK"

# Update the market / product code cells for the new UK sheet.
$uk.Range("B2").Value = "UK Market"
$uk.Range("B4").Value = "NGC-2741/T3356/T3357"

# Leave the source "Poland" sheet with a "select all" selection state (as it
# was right before being copied), and place the cursor on B4 / make "UK" the
# active sheet and tab.
$poland.Cells.Select()

$uk.Activate()
$uk.Range("B4").Select()
